$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update threshold values
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 12
$ws.Range("C5").Value = 18

# Update the active selection to E4
$ws.Range("E4").Select()

$wb.Save()
